$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 00:16"

# Country label swaps (shared-string reordering in the diff manifests as the
# country name at these two rows trading places, while the numeric columns
# for the same two rows follow the pattern of new data landing on the first
# row and old data shifting down into the second row).
$ws.Range("A51").Value = "Barein"
$ws.Range("A52").Value = "Honduras"

$ws.Range("A81").Value = "Bulgaria"
$ws.Range("A82").Value = "Bosnia y Herzegovina"

$ws.Range("A129").Value = "Ruanda"
$ws.Range("A130").Value = "Namibia"

# Updated numeric statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected rows

# Row 4
$ws.Range("B4").Value = 4488483
$ws.Range("C4").Value = 55073
$ws.Range("D4").Value = 2172095
$ws.Range("E4").Value = 2164292
$ws.Range("G4").Value = 1021
$ws.Range("H4").Value = 152096

# Row 5
$ws.Range("B5").Value = 2480888
$ws.Range("C5").Value = 37408
$ws.Range("E5").Value = 724751
$ws.Range("G5").Value = 791
$ws.Range("H5").Value = 88470

# Row 21
$ws.Range("B21").Value = 207951
$ws.Range("C21").Value = 572
$ws.Range("E21").Value = 7344

# Row 28
$ws.Range("B28").Value = 92947
$ws.Range("C28").Value = 465
$ws.Range("D28").Value = 35959
$ws.Range("E28").Value = 52297
$ws.Range("G28").Value = 39
$ws.Range("H28").Value = 4691

# Row 51
$ws.Range("B51").Value = 39921
$ws.Range("C51").Value = 439
$ws.Range("D51").Value = 36531
$ws.Range("E51").Value = 3249
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 141

# Row 52
$ws.Range("B52").Value = 39741
$ws.Range("C52").Value = 465
$ws.Range("D52").Value = 5039
$ws.Range("E52").Value = 33536
$ws.Range("G52").Value = 50
$ws.Range("H52").Value = 1166

# Row 55
$ws.Range("D55").Value = 31000
$ws.Range("E55").Value = 1631

# Row 81
$ws.Range("B81").Value = 10871
$ws.Range("C81").Value = 250
$ws.Range("D81").Value = 5766
$ws.Range("E81").Value = 4750
$ws.Range("G81").Value = 8
$ws.Range("H81").Value = 355

# Row 82
$ws.Range("B82").Value = 10766
$ws.Range("C82").Value = 268
$ws.Range("D82").Value = 5220
$ws.Range("E82").Value = 5249
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 297

# Row 129
$ws.Range("B129").Value = 1926
$ws.Range("C129").Value = 47
$ws.Range("D129").Value = 1005
$ws.Range("E129").Value = 916
$ws.Range("H129").Value = 5

# Row 130
$ws.Range("B130").Value = 1917
$ws.Range("C130").Value = 74
$ws.Range("D130").Value = 104
$ws.Range("E130").Value = 1805
$ws.Range("H130").Value = 8
